$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "1.002", "250.90") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.348.54"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.932.93"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "250.90"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").Value = "0.7155"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "0.3268"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "27.61"
$ws.Range("E9").Value = "  +5.22%  "
$ws.Range("D10").Value = "0.07173"
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("D11").Value = "0.8007"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").Value = "0.08077"
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").Value = "1.933.37"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "94.54"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "14.88"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("D17").Value = "30.341.51"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "252.36"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").Value = "0.000008113"
$ws.Range("E19").Value = "  +3.28%  "
$ws.Range("D20").Value = "5.809"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "2.185.97"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "6.928"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").Value = "9.702"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").Value = "165.44"
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("E27").Value = "  +5.08%  "
$ws.Range("D28").Value = "19.20"
$ws.Range("E28").Value = "  +2.76%  "
$ws.Range("D29").Value = "0.1288"
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("D30").Value = "1.366"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").Value = "1.542"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "4.421"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "4.197"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("D34").Value = "0.05197"
$ws.Range("E34").Value = "  +3.61%  "
$ws.Range("D35").Value = "1.270"
$ws.Range("E35").Value = "  +6.91%  "
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("D37").Value = "2.764"
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("D38").Value = "0.01960"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "2.799"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "79.02"
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").Value = "0.4521"
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("D43").Value = "2.026"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "0.8398"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").Value = "101.77"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "9.760"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").Value = "7.411"
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.06063"
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.4176"
$ws.Range("E51").Value = "  +2.57%  "

# Remove the temporary text-number-format so column D keeps its
# original (default) style, matching the source workbook.
$ws.Range("D2:D51").ClearFormats()

